# Scheduled-runner refresh of cached market-price / profit columns
# (H, I, J, K, L, M, N) across the per-job leve sheets. Values below are
# the new market snapshot; row identity (A-G) is untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1921.3572
$ws.Range("J19").Value = 3416.5
$ws.Range("L19").Value = 3416.5
$ws.Range("N19").Value = -3766.5

$ws.Range("H51").Value = 7455.6665
$ws.Range("J51").Value = 8285.444
$ws.Range("L51").Value = 8285.444
$ws.Range("N51").Value = -9253.444

$ws.Range("H101").Value = 666.7273
$ws.Range("J101").Value = 601
$ws.Range("L101").Value = 1803
$ws.Range("N101").Value = -5047

$ws.Range("H132").Value = 2554.92
$ws.Range("I132").Value = 2554.92
$ws.Range("K132").Value = 7664.76
$ws.Range("M132").Value = -5134.76

$ws.Range("H133").Value = 82237.75
$ws.Range("J133").Value = 89983.664
$ws.Range("L133").Value = 89983.664
$ws.Range("N133").Value = -100103.664

$ws.Range("H137").Value = 2762.6897
$ws.Range("I137").Value = 1814.9375
$ws.Range("K137").Value = 5444.8125
$ws.Range("M137").Value = -2894.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 309.5
$ws.Range("I4").Value = 309.5
$ws.Range("K4").Value = 309.5
$ws.Range("M4").Value = -193.5

$ws.Range("H32").Value = 6339.828
$ws.Range("I32").Value = 2707.9153
$ws.Range("K32").Value = 2707.9153
$ws.Range("M32").Value = -2420.9153

$ws.Range("H45").Value = 1784.7693
$ws.Range("I45").Value = 1114.619
$ws.Range("J45").Value = 4599.4
$ws.Range("K45").Value = 1114.619
$ws.Range("L45").Value = 4599.4
$ws.Range("M45").Value = -737.6189999999999
$ws.Range("N45").Value = -5353.4

$ws.Range("H61").Value = 6748.9624
$ws.Range("I61").Value = 5898.1396
$ws.Range("K61").Value = 5898.1396
$ws.Range("M61").Value = -5686.1396

$ws.Range("H74").Value = 3755.8462
$ws.Range("I74").Value = 1769.6666
$ws.Range("K74").Value = 1769.6666
$ws.Range("M74").Value = -895.6666

$ws.Range("H77").Value = 3755.8462
$ws.Range("I77").Value = 1769.6666
$ws.Range("K77").Value = 8848.333000000001
$ws.Range("M77").Value = -4480.333000000001

$ws.Range("H120").Value = 199950
$ws.Range("J120").Value = 199950
$ws.Range("L120").Value = 199950
$ws.Range("N120").Value = -209626

$ws.Range("H122").Value = 3637.5186
$ws.Range("I122").Value = 3473.6667
$ws.Range("K122").Value = 10421.0001
$ws.Range("M122").Value = -7971.000100000001

$ws.Range("H132").Value = 3337.2424
$ws.Range("I132").Value = 1934.238
$ws.Range("J132").Value = 5792.5
$ws.Range("K132").Value = 5802.714
$ws.Range("L132").Value = 17377.5
$ws.Range("M132").Value = -3272.714
$ws.Range("N132").Value = -22437.5

$ws.Range("H136").Value = 6748.9624
$ws.Range("I136").Value = 5898.1396
$ws.Range("K136").Value = 17694.4188
$ws.Range("M136").Value = -15144.4188

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 125000
$ws.Range("J33").Value = 125000
$ws.Range("L33").Value = 125000
$ws.Range("N33").Value = -125672

$ws.Range("H61").Value = 49999
$ws.Range("J61").Value = 49999
$ws.Range("L61").Value = 49999
$ws.Range("N61").Value = -50625

$ws.Range("H86").Value = 1239.909
$ws.Range("I86").Value = 1153.5625
$ws.Range("K86").Value = 1153.5625
$ws.Range("M86").Value = -30.5625

$ws.Range("H87").Value = 199950
$ws.Range("J87").Value = 199950
$ws.Range("L87").Value = 199950
$ws.Range("N87").Value = -202446

$ws.Range("H89").Value = 1239.909
$ws.Range("I89").Value = 1153.5625
$ws.Range("K89").Value = 5767.8125
$ws.Range("M89").Value = -151.8125

$ws.Range("H90").Value = 199950
$ws.Range("J90").Value = 199950
$ws.Range("L90").Value = 599850
$ws.Range("N90").Value = -612330

$ws.Range("H97").Value = 137344.33
$ws.Range("I97").Value = 12133.333
$ws.Range("J97").Value = 199949.83
$ws.Range("K97").Value = 12133.333
$ws.Range("L97").Value = 199949.83
$ws.Range("M97").Value = -11142.333
$ws.Range("N97").Value = -201931.83

$ws.Range("H100").Value = 188398.67
$ws.Range("J100").Value = 188398.67
$ws.Range("L100").Value = 188398.67
$ws.Range("N100").Value = -190562.67

$ws.Range("H103").Value = 108943.1
$ws.Range("J103").Value = 108943.1
$ws.Range("L103").Value = 108943.1
$ws.Range("N103").Value = -111287.1

$ws.Range("H111").Value = 199950
$ws.Range("J111").Value = 199950
$ws.Range("L111").Value = 199950
$ws.Range("N111").Value = -208130

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 621.5
$ws.Range("I7").Value = 707.4783
$ws.Range("J7").Value = 441.72726
$ws.Range("K7").Value = 707.4783
$ws.Range("L7").Value = 441.72726
$ws.Range("M7").Value = -594.4783
$ws.Range("N7").Value = -667.72726

$ws.Range("H58").Value = 6429.7407
$ws.Range("I58").Value = 2674
$ws.Range("K58").Value = 2674
$ws.Range("M58").Value = -2471

$ws.Range("H99").Value = 2492.8667
$ws.Range("I99").Value = 1990.3636
$ws.Range("K99").Value = 1990.3636
$ws.Range("M99").Value = -492.3635999999999

$ws.Range("H105").Value = 2422.8333
$ws.Range("J105").Value = 2222
$ws.Range("L105").Value = 2222
$ws.Range("N105").Value = -5716

$ws.Range("H126").Value = 2492.8667
$ws.Range("I126").Value = 1990.3636
$ws.Range("K126").Value = 5971.0908
$ws.Range("M126").Value = -3501.0908

$ws.Range("H136").Value = 6429.7407
$ws.Range("I136").Value = 2674
$ws.Range("K136").Value = 8022
$ws.Range("M136").Value = -5472

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I112").Value = 13979.857
$ws.Range("J112").Value = 20000
$ws.Range("K112").Value = 41939.571
$ws.Range("L112").Value = 60000
$ws.Range("M112").Value = -40831.571
$ws.Range("N112").Value = -62216

$ws.Range("H113").Value = 1633.5555
$ws.Range("J113").Value = 1665.375
$ws.Range("L113").Value = 4996.125
$ws.Range("N113").Value = -9336.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6650.7188
$ws.Range("J80").Value = 6821.1
$ws.Range("L80").Value = 6821.1
$ws.Range("N80").Value = -8817.1

$ws.Range("H83").Value = 6650.7188
$ws.Range("J83").Value = 6821.1
$ws.Range("L83").Value = 34105.5
$ws.Range("N83").Value = -44089.5

$ws.Range("H113").Value = 4212.5625
$ws.Range("I113").Value = 3727.6667
$ws.Range("J113").Value = 5667.25
$ws.Range("K113").Value = 3727.6667
$ws.Range("L113").Value = 5667.25
$ws.Range("M113").Value = -1557.6667
$ws.Range("N113").Value = -10007.25

$ws.Range("H123").Value = 64750
$ws.Range("J123").Value = 64750
$ws.Range("L123").Value = 64750
$ws.Range("N123").Value = -69650

$ws.Range("H126").Value = 4294.5
$ws.Range("J126").Value = 7869
$ws.Range("L126").Value = 23607
$ws.Range("N126").Value = -28547

$ws.Range("H132").Value = 5494.625
$ws.Range("I132").Value = 5017.3213
$ws.Range("K132").Value = 15051.9639
$ws.Range("M132").Value = -12521.9639

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1995
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H22").Value = 1362.9445
$ws.Range("I22").Value = 640.8182
$ws.Range("J22").Value = 2497.7144
$ws.Range("K22").Value = 640.8182
$ws.Range("L22").Value = 2497.7144
$ws.Range("M22").Value = -345.8182
$ws.Range("N22").Value = -3087.7144

$ws.Range("H27").Value = 1362.9445
$ws.Range("I27").Value = 640.8182
$ws.Range("J27").Value = 2497.7144
$ws.Range("K27").Value = 640.8182
$ws.Range("L27").Value = 2497.7144
$ws.Range("M27").Value = -533.8182
$ws.Range("N27").Value = -2711.7144

$ws.Range("H46").Value = 3407.6924
$ws.Range("I46").Value = 1540
$ws.Range("K46").Value = 1540
$ws.Range("M46").Value = -1352

$ws.Range("H93").Value = 1659
$ws.Range("J93").Value = 2240.4443
$ws.Range("L93").Value = 2240.4443
$ws.Range("N93").Value = -4736.4443

$ws.Range("H132").Value = 5813.84
$ws.Range("I132").Value = 5306.0835
$ws.Range("K132").Value = 15918.2505
$ws.Range("M132").Value = -13388.2505

$ws.Range("H136").Value = 5378.854
$ws.Range("I136").Value = 4120.921
$ws.Range("K136").Value = 12362.763
$ws.Range("M136").Value = -9812.763000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 17283.166
$ws.Range("J41").Value = 17319.8
$ws.Range("L41").Value = 17319.8
$ws.Range("N41").Value = -18099.8

$ws.Range("H107").Value = 1251750.5
$ws.Range("I107").Value = 2224141
$ws.Range("K107").Value = 6672423
$ws.Range("M107").Value = -6670503

$ws.Range("H113").Value = 935.44446
$ws.Range("I113").Value = 511.33334
$ws.Range("J113").Value = 1783.6666
$ws.Range("K113").Value = 1534.00002
$ws.Range("L113").Value = 5350.9998
$ws.Range("M113").Value = 635.9999800000001
$ws.Range("N113").Value = -9690.9998
